$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4677.615
$ws.Range("I98").Value = 2580.9
$ws.Range("J98").Value = 11666.667
$ws.Range("K98").Value = 2580.9
$ws.Range("L98").Value = 11666.667
$ws.Range("M98").Value = -1082.9
$ws.Range("N98").Value = -14662.667
$ws.Range("H122").Value = 4677.615
$ws.Range("I122").Value = 2580.9
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 7742.700000000001
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -5292.700000000001
$ws.Range("N122").Value = -39900.001
$ws.Range("H135").Value = 817.48486
$ws.Range("I135").Value = 569.11536
$ws.Range("J135").Value = 1740
$ws.Range("K135").Value = 5122.03824
$ws.Range("L135").Value = 15660
$ws.Range("M135").Value = -2587.03824
$ws.Range("N135").Value = -20730
$ws.Range("H137").Value = 2816.814
$ws.Range("I137").Value = 2047.12
$ws.Range("J137").Value = 3885.8333
$ws.Range("K137").Value = 6141.36
$ws.Range("L137").Value = 11657.4999
$ws.Range("M137").Value = -3591.36
$ws.Range("N137").Value = -16757.4999
$ws.Range("H138").Value = 4066.4824
$ws.Range("I138").Value = 3750.8572
$ws.Range("J138").Value = 4094.8076
$ws.Range("K138").Value = 11252.5716
$ws.Range("L138").Value = 12284.4228
$ws.Range("M138").Value = -6112.571599999999
$ws.Range("N138").Value = -22564.4228
$ws.Range("H141").Value = 1947.2667
$ws.Range("I141").Value = 1689.1538
$ws.Range("J141").Value = 3625
$ws.Range("K141").Value = 5067.4614
$ws.Range("L141").Value = 10875
$ws.Range("M141").Value = 112.5385999999999
$ws.Range("N141").Value = -21235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8774.84
$ws.Range("I32").Value = 4859.323
$ws.Range("J32").Value = 16046.515
$ws.Range("K32").Value = 4859.323
$ws.Range("L32").Value = 16046.515
$ws.Range("M32").Value = -4572.323
$ws.Range("N32").Value = -16620.515
$ws.Range("H45").Value = 1100.3334
$ws.Range("I45").Value = 848.6667
$ws.Range("K45").Value = 848.6667
$ws.Range("M45").Value = -471.6667
$ws.Range("H63").Value = 19789602
$ws.Range("I63").Value = 23087286
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 23087286
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -23086600
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 19789602
$ws.Range("I66").Value = 23087286
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 115436430
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -115432998
$ws.Range("N66").Value = -24364
$ws.Range("H132").Value = 2382.3572
$ws.Range("I132").Value = 1227.8572
$ws.Range("J132").Value = 4691.357
$ws.Range("K132").Value = 3683.5716
$ws.Range("L132").Value = 14074.071
$ws.Range("M132").Value = -1153.5716
$ws.Range("N132").Value = -19134.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 454
$ws.Range("I22").Value = 306.66666
$ws.Range("K22").Value = 306.66666
$ws.Range("M22").Value = -133.66666
$ws.Range("H31").Value = 5031.5
$ws.Range("I31").Value = 1100
$ws.Range("J31").Value = 6342
$ws.Range("K31").Value = 1100
$ws.Range("L31").Value = 6342
$ws.Range("M31").Value = -848
$ws.Range("N31").Value = -6846
$ws.Range("H86").Value = 2933.3333
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2933.3333
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H94").Value = 22729086
$ws.Range("I94").Value = 38463340
$ws.Range("J94").Value = 1832.3334
$ws.Range("K94").Value = 38463340
$ws.Range("L94").Value = 1832.3334
$ws.Range("M94").Value = -38462889
$ws.Range("N94").Value = -2734.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4121.1763
$ws.Range("I31").Value = 1588.3334
$ws.Range("K31").Value = 1588.3334
$ws.Range("M31").Value = -1293.3334
$ws.Range("H34").Value = 4121.1763
$ws.Range("I34").Value = 1588.3334
$ws.Range("K34").Value = 1588.3334
$ws.Range("M34").Value = -1386.3334
$ws.Range("H39").Value = 17131.234
$ws.Range("J39").Value = 24852.908
$ws.Range("L39").Value = 24852.908
$ws.Range("N39").Value = -25634.908
$ws.Range("H49").Value = 17131.234
$ws.Range("J49").Value = 24852.908
$ws.Range("L49").Value = 24852.908
$ws.Range("N49").Value = -25216.908
$ws.Range("H58").Value = 2191.4167
$ws.Range("I58").Value = 1745.3684
$ws.Range("J58").Value = 10666.333
$ws.Range("K58").Value = 1745.3684
$ws.Range("L58").Value = 10666.333
$ws.Range("M58").Value = -1542.3684
$ws.Range("N58").Value = -11072.333
$ws.Range("H68").Value = 46526.9
$ws.Range("J68").Value = 46526.9
$ws.Range("L68").Value = 46526.9
$ws.Range("N68").Value = -48024.9
$ws.Range("H71").Value = 46526.9
$ws.Range("J71").Value = 46526.9
$ws.Range("L71").Value = 139580.7
$ws.Range("N71").Value = -147068.7
$ws.Range("H132").Value = 3913.5667
$ws.Range("I132").Value = 3685.7368
$ws.Range("J132").Value = 4307.091
$ws.Range("K132").Value = 11057.2104
$ws.Range("L132").Value = 12921.273
$ws.Range("M132").Value = -8527.2104
$ws.Range("N132").Value = -17981.273
$ws.Range("H136").Value = 2191.4167
$ws.Range("I136").Value = 1745.3684
$ws.Range("J136").Value = 10666.333
$ws.Range("K136").Value = 5236.1052
$ws.Range("L136").Value = 31998.999
$ws.Range("M136").Value = -2686.1052
$ws.Range("N136").Value = -37098.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 10325.5
$ws.Range("I86").Value = 13633.333
$ws.Range("J86").Value = 7017.6665
$ws.Range("K86").Value = 40899.999
$ws.Range("L86").Value = 21052.9995
$ws.Range("M86").Value = -39713.999
$ws.Range("N86").Value = -23424.9995
$ws.Range("H89").Value = 10325.5
$ws.Range("I89").Value = 13633.333
$ws.Range("J89").Value = 7017.6665
$ws.Range("K89").Value = 122699.997
$ws.Range("L89").Value = 63158.9985
$ws.Range("M89").Value = -116771.997
$ws.Range("N89").Value = -75014.9985
$ws.Range("H113").Value = 598.5405
$ws.Range("I113").Value = 624.05884
$ws.Range("K113").Value = 1872.17652
$ws.Range("M113").Value = 297.82348

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2933.9546
$ws.Range("I122").Value = 2066.8333
$ws.Range("J122").Value = 6836
$ws.Range("K122").Value = 6200.499899999999
$ws.Range("L122").Value = 20508
$ws.Range("M122").Value = -3750.499899999999
$ws.Range("N122").Value = -25408

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6123.9414
$ws.Range("I7").Value = 3009.818
$ws.Range("J7").Value = 11833.167
$ws.Range("K7").Value = 3009.818
$ws.Range("L7").Value = 11833.167
$ws.Range("M7").Value = -2897.818
$ws.Range("N7").Value = -12057.167
$ws.Range("H10").Value = 37666.668
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 37666.668
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 37666.668
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -37946.668
$ws.Range("H68").Value = 922.06976
$ws.Range("I68").Value = 735.10254
$ws.Range("J68").Value = 2745
$ws.Range("K68").Value = 735.10254
$ws.Range("L68").Value = 2745
$ws.Range("M68").Value = 13.89746000000002
$ws.Range("N68").Value = -4243
$ws.Range("H71").Value = 922.06976
$ws.Range("I71").Value = 735.10254
$ws.Range("J71").Value = 2745
$ws.Range("K71").Value = 3675.5127
$ws.Range("L71").Value = 13725
$ws.Range("M71").Value = 68.48730000000023
$ws.Range("N71").Value = -21213
$ws.Range("H82").Value = 4826.231
$ws.Range("I82").Value = 5694.95
$ws.Range("J82").Value = 1930.5
$ws.Range("K82").Value = 5694.95
$ws.Range("L82").Value = 1930.5
$ws.Range("M82").Value = -5333.95
$ws.Range("N82").Value = -2652.5
$ws.Range("H85").Value = 4826.231
$ws.Range("I85").Value = 5694.95
$ws.Range("J85").Value = 1930.5
$ws.Range("K85").Value = 5694.95
$ws.Range("L85").Value = 1930.5
$ws.Range("M85").Value = -4446.95
$ws.Range("N85").Value = -4426.5
$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 5000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -7496
$ws.Range("H122").Value = 2932.606
$ws.Range("I122").Value = 2216.2917
$ws.Range("J122").Value = 4842.778
$ws.Range("K122").Value = 6648.875100000001
$ws.Range("L122").Value = 14528.334
$ws.Range("M122").Value = -4198.875100000001
$ws.Range("N122").Value = -19428.334
$ws.Range("H126").Value = 6123.9414
$ws.Range("I126").Value = 3009.818
$ws.Range("J126").Value = 11833.167
$ws.Range("K126").Value = 9029.454000000002
$ws.Range("L126").Value = 35499.501
$ws.Range("M126").Value = -6559.454000000002
$ws.Range("N126").Value = -40439.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1752.5
$ws.Range("I13").Value = 1005
$ws.Range("J13").Value = 2500
$ws.Range("K13").Value = 1005
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = -865
$ws.Range("N13").Value = -2780
$ws.Range("H81").Value = 53572660
$ws.Range("I81").Value = 53572660
$ws.Range("K81").Value = 107145320
$ws.Range("M81").Value = -107144259
$ws.Range("H84").Value = 53572660
$ws.Range("I84").Value = 53572660
$ws.Range("K84").Value = 535726600
$ws.Range("M84").Value = -535721296
$ws.Range("H92").Value = 39750
$ws.Range("J92").Value = 39750
$ws.Range("L92").Value = 39750
$ws.Range("N92").Value = -44742
$ws.Range("H93").Value = 39750
$ws.Range("J93").Value = 39750
$ws.Range("L93").Value = 39750
$ws.Range("N93").Value = -44742
$ws.Range("H109").Value = 28377
$ws.Range("J109").Value = 28377
$ws.Range("L109").Value = 28377
$ws.Range("N109").Value = -31151
$ws.Range("H122").Value = 4639.316
$ws.Range("I122").Value = 2677
$ws.Range("J122").Value = 7337.5
$ws.Range("K122").Value = 8031
$ws.Range("L122").Value = 22012.5
$ws.Range("M122").Value = -5581
$ws.Range("N122").Value = -26912.5
